$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "'36.48"
$ws.Range("C12").Value = "'11.28"
$ws.Range("D12").Value = "'47.76"

$ws.Range("B13").Value = "'6.83"
$ws.Range("C13").Value = "'39.18"
$ws.Range("D13").Value = "'46.01"

$ws.Range("B15").Value = "'75.24"
$ws.Range("C15").Value = "'23.26"
